# Apply the "update with new logo and colors" metadata revision to the
# CIBMTR Priority Variables (FY21) ValueSet workbook.
#
# Summary of the change (per the target diff):
#   - Metadata sheet: bump Version, change Status, bump Date, replace the
#     Contact text, add a second Contact row, add a new Jurisdiction row
#     (pushing Description/Purpose/Copyright/Immutable down by one row).
#   - "Include from LOINC" sheet: unaffected in content (only shared
#     string table bookkeeping shifts because of the Metadata edits).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Make room: shift the trailing rows (old 12-15) down to (13-16) ---
# Work bottom-up so we never overwrite a value before it has been copied.

# Old row 15 (Immutable / BooleanType[null]) -> new row 16 (brand new row,
# so copy formatting from the row above it first).
$ws.Range("A16").Value2 = $ws.Range("A15").Value2
$ws.Range("B16").Value2 = $ws.Range("B15").Value2
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

# Old row 14 (Copyright / empty) -> new row 15
$ws.Range("A15").Value2 = $ws.Range("A14").Value2
$ws.Range("B15").Value2 = $ws.Range("B14").Value2

# Old row 13 (Purpose / empty) -> new row 14
$ws.Range("A14").Value2 = $ws.Range("A13").Value2
$ws.Range("B14").Value2 = $ws.Range("B13").Value2

# Old row 12 (Description / Priority Variables for CIBMTR (FY21)) -> new row 13
$ws.Range("A13").Value2 = $ws.Range("A12").Value2
$ws.Range("B13").Value2 = $ws.Range("B12").Value2

# --- New row 12: Jurisdiction, with no value (matches "No display for ..." blank) ---
$ws.Range("A12").Value2 = "Jurisdiction"
$ws.Range("B12").Value2 = ""

# --- Row 11 already duplicated "Contact" -- turn it into the second contact line ---
$ws.Range("B11").Value2 = "Bob Milius (bmilius@nmdp.org)"

# --- Row 10: first Contact line gets the real ContactDetail display text ---
$ws.Range("B10").Value2 = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- Simple scalar metadata updates ---
$ws.Range("B3").Value2 = "0.1.7"
$ws.Range("B6").Value2 = "draft"
$ws.Range("B8").Value2 = "2024-08-23T10:17:11-05:00"
